$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("1碑影迷踪")
$ws.Range("B5").Value = "从燕墩所在位置向北望去半里左右会看到一座位于中轴线上的古代建筑"
